$d = $word.ActiveDocument

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    # Exclude the trailing cell-end mark from the range before mutating.
    $scoped = $d.Range($r.Start, $r.End - 1)
    $scoped.Text = $newText
}

# --- Header date paragraph ---
$d.Content.Find.Execute("2024-08-14 Wednesday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-08-15 Thursday", 2) | Out-Null

# --- Table of division problems ---
$t = $d.Tables.Item(1)

# Row 1
Set-CellText $t 1 1 "74÷8="
Set-CellText $t 1 2 "26÷5="
Set-CellText $t 1 3 "35÷2="
Set-CellText $t 1 4 "98÷8="
Set-CellText $t 1 5 "96÷6="

# Row 5
Set-CellText $t 5 1 "69÷4="
Set-CellText $t 5 2 "32÷2="
Set-CellText $t 5 3 "29÷9="
Set-CellText $t 5 4 "84÷7="
Set-CellText $t 5 5 "50÷6="

# Row 9
Set-CellText $t 9 1 "65÷2="
Set-CellText $t 9 2 "57÷7="
Set-CellText $t 9 3 "44÷4="
Set-CellText $t 9 4 "68÷4="
Set-CellText $t 9 5 "82÷8="

# Row 13
Set-CellText $t 13 1 "41÷4="
Set-CellText $t 13 2 "53÷6="
Set-CellText $t 13 3 "70÷5="
Set-CellText $t 13 4 "94÷2="
Set-CellText $t 13 5 "89÷4="

# Row 17
Set-CellText $t 17 1 "19÷4="
Set-CellText $t 17 2 "77÷3="
Set-CellText $t 17 3 "44÷2="
Set-CellText $t 17 4 "23÷6="
Set-CellText $t 17 5 "59÷3="
